$d = $word.ActiveDocument

# -----------------------------------------------------------------
# Change 1: the empty paragraph after "...app_main" (right before
# "But I don\t Know why I'm getting error") currently only holds the
# hidden "_GoBack" bookmark. Remove that bookmark so the paragraph
# becomes a truly empty paragraph.
# -----------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# -----------------------------------------------------------------
# Change 2: the "Start idf.py menuconfig" paragraph gets a hanging
# indent: left indent 1440 twips (72 pt) and first-line indent 720
# twips (36 pt).
# -----------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -eq "Start idf.py menuconfig`r") {
        $para.Range.ParagraphFormat.LeftIndent = 72
        $para.Range.ParagraphFormat.FirstLineIndent = 36
        break
    }
}

# -----------------------------------------------------------------
# Change 3: the "_GoBack" bookmark re-appears in the middle of the
# "For example you can use it to change cpu frequency of ESP32"
# paragraph, splitting "can" into "ca" | bookmark | "n".
# -----------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -like "For example you can use it to change cpu frequency of ESP32*") {
        $paraStart = $para.Range.Start
        $splitPos = $paraStart + ("For example you ca").Length

        # Insert the (hidden) bookmark right between "ca" and "n".
        $bmRange = $d.Range($splitPos, $splitPos)
        $d.Bookmarks.Add("_GoBack", $bmRange)

        # Force the now-split leading run ("For example you ca") to be
        # re-written so its <w:t> element drops the (no-longer-needed)
        # xml:space="preserve" flag, matching a genuine Word re-save.
        $leadRange = $d.Range($paraStart, $splitPos)
        $leadRange.Text = "#"
        $leadRange2 = $d.Range($paraStart, $paraStart + 1)
        $leadRange2.Text = "For example you ca"
        break
    }
}
